$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ "G"=[double]"0.2753413333333334"; "H"=[double]"0.8260240000000001"; "I"=[double]"0.006630378892106956"; "J"=[double]"0.006630378892106955"; "O"=[double]"0.004953627445374112"; "P"=[double]"0.004953627445374111"; "Q"=[double]"0.0009779206355555556"; "R"=[double]"0.008801285720000001"; "S"=[double]"3.284442685317021E-05"; "T"=[double]"3.284442685317021E-05" }
    3 = @{ "G"=[double]"0.2753413333333334"; "H"=[double]"0.8260240000000001"; "I"=[double]"0.006630378892106956"; "J"=[double]"0.006630378892106955"; "O"=[double]"0.0778289024983856"; "P"=[double]"0.07782890249838559"; "Q"=[double]"0.01536459708266667"; "R"=[double]"0.138281373744"; "S"=[double]"0.0005160351123211463"; "T"=[double]"0.0005160351123211461" }
    4 = @{ "G"=[double]"0.2753413333333334"; "H"=[double]"0.8260240000000001"; "I"=[double]"0.006630378892106956"; "J"=[double]"0.006630378892106955"; "M"=[double]"0.509521"; "N"=[double]"1.528563"; "O"=[double]"0.7106458591068409"; "P"=[double]"0.7106458591068406"; "Q"=[double]"0.1402921915013333"; "R"=[double]"1.262629723512"; "S"=[double]"0.004711851303985212"; "T"=[double]"0.00471185130398521" }
    5 = @{ "G"=[double]"0.2753413333333334"; "H"=[double]"0.8260240000000001"; "I"=[double]"0.006630378892106956"; "J"=[double]"0.006630378892106955"; "M"=[double]"0.06794533333333334"; "N"=[double]"0.203836"; "O"=[double]"0.09476561275976328"; "P"=[double]"0.09476561275976325"; "Q"=[double]"0.01870815867377778"; "R"=[double]"0.168373428064"; "S"=[double]"0.0006283319185399161"; "T"=[double]"0.0006283319185399158" }
    6 = @{ "G"=[double]"0.2753413333333334"; "H"=[double]"0.8260240000000001"; "I"=[double]"0.006630378892106956"; "J"=[double]"0.006630378892106955"; "M"=[double]"0.080163"; "N"=[double]"0.240489"; "O"=[double]"0.1118059981896363"; "P"=[double]"0.1118059981896363"; "Q"=[double]"0.022072187304"; "R"=[double]"0.198649685736"; "S"=[double]"0.0007413161304075131"; "T"=[double]"0.0007413161304075128" }
    7 = @{ "I"=[double]"0.03952244389885164"; "J"=[double]"0.03952244389885164"; "O"=[double]"0.004953627445374112"; "P"=[double]"0.004953627445374111"; "S"=[double]"0.0001957794628056101"; "T"=[double]"0.0001957794628056101" }
    8 = @{ "I"=[double]"0.03952244389885164"; "J"=[double]"0.03952244389885164"; "O"=[double]"0.0778289024983856"; "P"=[double]"0.07782890249838559"; "S"=[double]"0.003075988432701639"; "T"=[double]"0.003075988432701639" }
    9 = @{ "I"=[double]"0.03952244389885164"; "J"=[double]"0.03952244389885164"; "M"=[double]"0.509521"; "N"=[double]"1.528563"; "O"=[double]"0.7106458591068409"; "P"=[double]"0.7106458591068406"; "Q"=[double]"0.8362554174179999"; "R"=[double]"7.526298756761999"; "S"=[double]"0.02808646109850135"; "T"=[double]"0.02808646109850134" }
    10 = @{ "I"=[double]"0.03952244389885164"; "J"=[double]"0.03952244389885164"; "M"=[double]"0.06794533333333334"; "N"=[double]"0.203836"; "O"=[double]"0.09476561275976328"; "P"=[double]"0.09476561275976325"; "Q"=[double]"0.111515821896"; "R"=[double]"1.003642397064"; "S"=[double]"0.003745368613838044"; "T"=[double]"0.003745368613838043" }
    11 = @{ "I"=[double]"0.03952244389885164"; "J"=[double]"0.03952244389885164"; "M"=[double]"0.080163"; "N"=[double]"0.240489"; "O"=[double]"0.1118059981896363"; "P"=[double]"0.1118059981896363"; "Q"=[double]"0.131568165054"; "R"=[double]"1.184113485486"; "S"=[double]"0.00441884629100501"; "T"=[double]"0.004418846291005008" }
    12 = @{ "G"=[double]"23.78768866666667"; "H"=[double]"71.363066"; "I"=[double]"0.5728213302306416"; "J"=[double]"0.5728213302306416"; "O"=[double]"0.004953627445374112"; "P"=[double]"0.004953627445374111"; "Q"=[double]"0.08448594091444445"; "R"=[double]"0.76037346823"; "S"=[double]"0.002837543462726214"; "T"=[double]"0.002837543462726213" }
    13 = @{ "G"=[double]"23.78768866666667"; "H"=[double]"71.363066"; "I"=[double]"0.5728213302306416"; "J"=[double]"0.5728213302306416"; "O"=[double]"0.0778289024983856"; "P"=[double]"0.07782890249838559"; "Q"=[double]"1.327400602977333"; "R"=[double]"11.946605426796"; "S"=[double]"0.04458205545951615"; "T"=[double]"0.04458205545951614" }
    14 = @{ "G"=[double]"23.78768866666667"; "H"=[double]"71.363066"; "I"=[double]"0.5728213302306416"; "J"=[double]"0.5728213302306416"; "M"=[double]"0.509521"; "N"=[double]"1.528563"; "O"=[double]"0.7106458591068409"; "P"=[double]"0.7106458591068406"; "Q"=[double]"12.12032691712867"; "R"=[double]"109.082942254158"; "S"=[double]"0.4070731063364777"; "T"=[double]"0.4070731063364776" }
    15 = @{ "G"=[double]"23.78768866666667"; "H"=[double]"71.363066"; "I"=[double]"0.5728213302306416"; "J"=[double]"0.5728213302306416"; "M"=[double]"0.06794533333333334"; "N"=[double]"0.203836"; "O"=[double]"0.09476561275976328"; "P"=[double]"0.09476561275976325"; "Q"=[double]"1.616262435686223"; "R"=[double]"14.546361921176"; "S"=[double]"0.05428376436116947"; "T"=[double]"0.05428376436116945" }
    16 = @{ "G"=[double]"23.78768866666667"; "H"=[double]"71.363066"; "I"=[double]"0.5728213302306416"; "J"=[double]"0.5728213302306416"; "M"=[double]"0.080163"; "N"=[double]"0.240489"; "O"=[double]"0.1118059981896363"; "P"=[double]"0.1118059981896363"; "Q"=[double]"1.906892486586"; "R"=[double]"17.162032379274"; "S"=[double]"0.06404486061075218"; "T"=[double]"0.06404486061075215" }
    17 = @{ "G"=[double]"0.5982033333333333"; "H"=[double]"1.79461"; "I"=[double]"0.01440508298011203"; "J"=[double]"0.01440508298011203"; "O"=[double]"0.004953627445374112"; "P"=[double]"0.004953627445374111"; "Q"=[double]"0.002124618838888889"; "R"=[double]"0.01912156955"; "S"=[double]"7.135741440317447E-05"; "T"=[double]"7.135741440317446E-05" }
    18 = @{ "G"=[double]"0.5982033333333333"; "H"=[double]"1.79461"; "I"=[double]"0.01440508298011203"; "J"=[double]"0.01440508298011203"; "O"=[double]"0.0778289024983856"; "P"=[double]"0.07782890249838559"; "Q"=[double]"0.03338094240666666"; "R"=[double]"0.30042848166"; "S"=[double]"0.001121131798740293"; "T"=[double]"0.001121131798740293" }
    19 = @{ "G"=[double]"0.5982033333333333"; "H"=[double]"1.79461"; "I"=[double]"0.01440508298011203"; "J"=[double]"0.01440508298011203"; "M"=[double]"0.509521"; "N"=[double]"1.528563"; "O"=[double]"0.7106458591068409"; "P"=[double]"0.7106458591068406"; "Q"=[double]"0.3047971606033333"; "R"=[double]"2.74317444543"; "S"=[double]"0.01023691256990705"; "T"=[double]"0.01023691256990704" }
    20 = @{ "G"=[double]"0.5982033333333333"; "H"=[double]"1.79461"; "I"=[double]"0.01440508298011203"; "J"=[double]"0.01440508298011203"; "M"=[double]"0.06794533333333334"; "N"=[double]"0.203836"; "O"=[double]"0.09476561275976328"; "P"=[double]"0.09476561275976325"; "Q"=[double]"0.04064512488444445"; "R"=[double]"0.3658061239600001"; "S"=[double]"0.001365106515465554"; "T"=[double]"0.001365106515465553" }
    21 = @{ "G"=[double]"0.5982033333333333"; "H"=[double]"1.79461"; "I"=[double]"0.01440508298011203"; "J"=[double]"0.01440508298011203"; "M"=[double]"0.080163"; "N"=[double]"0.240489"; "O"=[double]"0.1118059981896363"; "P"=[double]"0.1118059981896363"; "Q"=[double]"0.04795377381"; "R"=[double]"0.43158396429"; "S"=[double]"0.001610574681595967"; "T"=[double]"0.001610574681595966" }
    22 = @{ "G"=[double]"15.22474833333333"; "H"=[double]"45.674245"; "I"=[double]"0.3666207639982877"; "J"=[double]"0.3666207639982877"; "O"=[double]"0.004953627445374112"; "P"=[double]"0.004953627445374111"; "Q"=[double]"0.05407323116388888"; "R"=[double]"0.4866590804749999"; "S"=[double]"0.001816102678585943"; "T"=[double]"0.001816102678585943" }
    23 = @{ "G"=[double]"15.22474833333333"; "H"=[double]"45.674245"; "I"=[double]"0.3666207639982877"; "J"=[double]"0.3666207639982877"; "O"=[double]"0.0778289024983856"; "P"=[double]"0.07782890249838559"; "Q"=[double]"0.8495714064966666"; "R"=[double]"7.64614265847"; "S"=[double]"0.02853369169510638"; "T"=[double]"0.02853369169510636" }
    24 = @{ "G"=[double]"15.22474833333333"; "H"=[double]"45.674245"; "I"=[double]"0.3666207639982877"; "J"=[double]"0.3666207639982877"; "M"=[double]"0.509521"; "N"=[double]"1.528563"; "O"=[double]"0.7106458591068409"; "P"=[double]"0.7106458591068406"; "Q"=[double]"7.757328995548333"; "R"=[double]"69.815960959935"; "S"=[double]"0.2605375277979696"; "T"=[double]"0.2605375277979694" }
    25 = @{ "G"=[double]"15.22474833333333"; "H"=[double]"45.674245"; "I"=[double]"0.3666207639982877"; "J"=[double]"0.3666207639982877"; "M"=[double]"0.06794533333333334"; "N"=[double]"0.203836"; "O"=[double]"0.09476561275976328"; "P"=[double]"0.09476561275976325"; "Q"=[double]"1.034450600424444"; "R"=[double]"9.31005540382"; "S"=[double]"0.0347430413507503"; "T"=[double]"0.03474304135075029" }
    26 = @{ "G"=[double]"15.22474833333333"; "H"=[double]"45.674245"; "I"=[double]"0.3666207639982877"; "J"=[double]"0.3666207639982877"; "M"=[double]"0.080163"; "N"=[double]"0.240489"; "O"=[double]"0.1118059981896363"; "P"=[double]"0.1118059981896363"; "Q"=[double]"1.220461500645"; "R"=[double]"10.984153505805"; "S"=[double]"0.06404486061075218"; "T"=[double]"0.06404486061075215" }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $updates[$row][$col]
    }
}

Write-Output "done"